$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-49 down to 28-50.
$ws.Rows.Item(27).Insert()

# Fill in the new row 27 with the new weekly record.
$ws.Range("A27").Value = 4
$ws.Range("B27").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C27").Value = "Los Lagos"
$ws.Range("D27").Value = 45233
$ws.Range("E27").Value = 10
$ws.Range("F27").Value = 100112013
$ws.Range("G27").Value = "Alcachofa"
$ws.Range("H27").Value = "Española"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 13000
$ws.Range("L27").Value = 13000
$ws.Range("M27").Value = 13000
$ws.Range("N27").Value = "$/caja 30 unidades"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 433
$ws.Range("Q27").Value = 30
$ws.Range("R27").Value = "Hortaliza"
